# Refactor & add CAN connection status
#
# The "FPS: <value>" and "0" translation rows (rows 29-30) are removed
# from the Translation sheet, shifting all subsequent rows up by two.
# The row that ends up last (row 35) is then repurposed to hold a brand
# new "CAN DISCONNECTED" text entry (with its own new Text ID and a
# "Small" typography instead of "Medium").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Remove the obsolete "FPS: <value>" and "0" rows (rows 29-30); this
# shifts rows 31-36 up to become rows 29-34.
$ws.Rows("29:30").Delete()

# Turn the new last row (35) into the new CAN-disconnected status text.
$ws.Range("B35").Value = "SingleUseId55"
$ws.Range("C35").Value = "Small"
$ws.Range("D35").Value = "Left"
$ws.Range("E35").Value = "LTR"
$ws.Range("F35").Value = "CAN DISCONNECTED"
